$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1825
$ws.Range("J19").Value = 2290.4
$ws.Range("L19").Value = 2290.4
$ws.Range("N19").Value = -2640.4
$ws.Range("H40").Value = 2912.4
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2825
$ws.Range("H52").Value = 3328.25
$ws.Range("I52").Value = 3328.25
$ws.Range("K52").Value = 9984.75
$ws.Range("M52").Value = -9824.75
$ws.Range("H125").Value = 1450
$ws.Range("J125").Value = 1395.1666
$ws.Range("L125").Value = 12556.4994
$ws.Range("N125").Value = -17476.4994
$ws.Range("H138").Value = 2456.2046
$ws.Range("I138").Value = 2455.577
$ws.Range("J138").Value = 2457.111
$ws.Range("K138").Value = 7366.731000000001
$ws.Range("L138").Value = 7371.333
$ws.Range("M138").Value = -2226.731000000001
$ws.Range("N138").Value = -17651.333
$ws.Range("H141").Value = 4635.1816
$ws.Range("I141").Value = 2847.1667
$ws.Range("K141").Value = 8541.500100000001
$ws.Range("M141").Value = -3361.500100000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2581
$ws.Range("I32").Value = 1762.1
$ws.Range("J32").Value = 9132.200000000001
$ws.Range("K32").Value = 1762.1
$ws.Range("L32").Value = 9132.200000000001
$ws.Range("M32").Value = -1475.1
$ws.Range("N32").Value = -9706.200000000001
$ws.Range("H45").Value = 1337.7307
$ws.Range("I45").Value = 1088.2778
$ws.Range("K45").Value = 1088.2778
$ws.Range("M45").Value = -711.2778000000001
$ws.Range("H88").Value = 4041.2727
$ws.Range("I88").Value = 1828
$ws.Range("J88").Value = 4533.1113
$ws.Range("K88").Value = 1828
$ws.Range("L88").Value = 4533.1113
$ws.Range("M88").Value = -1422
$ws.Range("N88").Value = -5345.1113
$ws.Range("H91").Value = 4041.2727
$ws.Range("I91").Value = 1828
$ws.Range("J91").Value = 4533.1113
$ws.Range("K91").Value = 1828
$ws.Range("L91").Value = 4533.1113
$ws.Range("M91").Value = -424
$ws.Range("N91").Value = -7341.1113
$ws.Range("H97").Value = 1053.8889
$ws.Range("I97").Value = 1053.8889
$ws.Range("K97").Value = 1053.8889
$ws.Range("M97").Value = -557.8888999999999
$ws.Range("H123").Value = 70500
$ws.Range("J123").Value = 70500
$ws.Range("L123").Value = 70500
$ws.Range("N123").Value = -80300
$ws.Range("H132").Value = 2131.923
$ws.Range("I132").Value = 1340.8334
$ws.Range("J132").Value = 3911.875
$ws.Range("K132").Value = 4022.5002
$ws.Range("L132").Value = 11735.625
$ws.Range("M132").Value = -1492.5002
$ws.Range("N132").Value = -16795.625

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1946.4166
$ws.Range("I20").Value = 1866.9445
$ws.Range("K20").Value = 1866.9445
$ws.Range("M20").Value = -1619.9445
$ws.Range("H86").Value = 128211.06
$ws.Range("I86").Value = 6137.6
$ws.Range("K86").Value = 6137.6
$ws.Range("M86").Value = -5014.6
$ws.Range("H89").Value = 128211.06
$ws.Range("I89").Value = 6137.6
$ws.Range("K89").Value = 30688
$ws.Range("M89").Value = -25072
$ws.Range("H94").Value = 1101.2
$ws.Range("I94").Value = 1126.5
$ws.Range("K94").Value = 1126.5
$ws.Range("M94").Value = -675.5
$ws.Range("H139").Value = 44999.5
$ws.Range("J139").Value = 44999.5
$ws.Range("L139").Value = 44999.5
$ws.Range("N139").Value = -55279.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2752.7917
$ws.Range("I31").Value = 2475.2942
$ws.Range("J31").Value = 3426.7144
$ws.Range("K31").Value = 2475.2942
$ws.Range("L31").Value = 3426.7144
$ws.Range("M31").Value = -2180.2942
$ws.Range("N31").Value = -4016.7144
$ws.Range("H34").Value = 2752.7917
$ws.Range("I34").Value = 2475.2942
$ws.Range("J34").Value = 3426.7144
$ws.Range("K34").Value = 2475.2942
$ws.Range("L34").Value = 3426.7144
$ws.Range("M34").Value = -2273.2942
$ws.Range("N34").Value = -3830.7144
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2105.125
$ws.Range("I132").Value = 1206.9
$ws.Range("K132").Value = 3620.7
$ws.Range("M132").Value = -1090.7

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 47
$ws.Range("I8").Value = 47
$ws.Range("K8").Value = 141
$ws.Range("M8").Value = -2
$ws.Range("H68").Value = 607.53845
$ws.Range("I68").Value = 574.8333
$ws.Range("K68").Value = 1724.4999
$ws.Range("M68").Value = -913.4999
$ws.Range("H71").Value = 607.53845
$ws.Range("I71").Value = 574.8333
$ws.Range("K71").Value = 5173.4997
$ws.Range("M71").Value = -1117.4997
$ws.Range("H107").Value = 871.5
$ws.Range("J107").Value = 1038.3846
$ws.Range("L107").Value = 3115.1538
$ws.Range("N107").Value = -6955.1538

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2812.389
$ws.Range("I7").Value = 1914
$ws.Range("K7").Value = 1914
$ws.Range("M7").Value = -1802
$ws.Range("H46").Value = 1849.1333
$ws.Range("I46").Value = 1211
$ws.Range("K46").Value = 1211
$ws.Range("M46").Value = -1023
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H126").Value = 2812.389
$ws.Range("I126").Value = 1914
$ws.Range("K126").Value = 5742
$ws.Range("M126").Value = -3272
$ws.Range("H136").Value = 4622.316
$ws.Range("I136").Value = 2991.8572
$ws.Range("K136").Value = 8975.571599999999
$ws.Range("M136").Value = -6425.571599999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("H49").Value = 70062
$ws.Range("J49").Value = 70062
$ws.Range("L49").Value = 70062
$ws.Range("N49").Value = -70522
$ws.Range("H70").Value = 25422
$ws.Range("J70").Value = 29777.5
$ws.Range("L70").Value = 29777.5
$ws.Range("N70").Value = -30407.5
$ws.Range("H73").Value = 25422
$ws.Range("J73").Value = 29777.5
$ws.Range("L73").Value = 29777.5
$ws.Range("N73").Value = -31961.5
$ws.Range("H122").Value = 113852.71
$ws.Range("I122").Value = 158644
$ws.Range("J122").Value = 1874.5
$ws.Range("K122").Value = 475932
$ws.Range("L122").Value = 5623.5
$ws.Range("M122").Value = -473482
$ws.Range("N122").Value = -10523.5
$ws.Range("H132").Value = 2632.6667
$ws.Range("I132").Value = 1449.5
$ws.Range("K132").Value = 4348.5
$ws.Range("M132").Value = -1818.5
$ws.Range("H136").Value = 11575354
$ws.Range("I136").Value = 15433109
$ws.Range("K136").Value = 46299327
$ws.Range("M136").Value = -46296777
$ws.Range("H139").Value = 67950
$ws.Range("J139").Value = 67950
$ws.Range("L139").Value = 67950
$ws.Range("N139").Value = -78230
